$d = $word.ActiveDocument

$d.Content.Find.Execute("872×6=", $true, $false, $false, $false, $false, $true, 1, $false, "709×6=", 2)
$d.Content.Find.Execute("496×6=", $true, $false, $false, $false, $false, $true, 1, $false, "556×6=", 2)
$d.Content.Find.Execute("811×5=", $true, $false, $false, $false, $false, $true, 1, $false, "733×5=", 2)
$d.Content.Find.Execute("368×6=", $true, $false, $false, $false, $false, $true, 1, $false, "258×3=", 2)
$d.Content.Find.Execute("159×2=", $true, $false, $false, $false, $false, $true, 1, $false, "464×2=", 2)
$d.Content.Find.Execute("236×4=", $true, $false, $false, $false, $false, $true, 1, $false, "559×7=", 2)
$d.Content.Find.Execute("421×5=", $true, $false, $false, $false, $false, $true, 1, $false, "234×3=", 2)
$d.Content.Find.Execute("252×6=", $true, $false, $false, $false, $false, $true, 1, $false, "706×7=", 2)
$d.Content.Find.Execute("471×7=", $true, $false, $false, $false, $false, $true, 1, $false, "749×2=", 2)
$d.Content.Find.Execute("148×6=", $true, $false, $false, $false, $false, $true, 1, $false, "695×5=", 2)
$d.Content.Find.Execute("234×2=", $true, $false, $false, $false, $false, $true, 1, $false, "668×4=", 2)
$d.Content.Find.Execute("718×2=", $true, $false, $false, $false, $false, $true, 1, $false, "745×3=", 2)
$d.Content.Find.Execute("848×2=", $true, $false, $false, $false, $false, $true, 1, $false, "347×2=", 2)
$d.Content.Find.Execute("163×3=", $true, $false, $false, $false, $false, $true, 1, $false, "444×7=", 2)
$d.Content.Find.Execute("894×6=", $true, $false, $false, $false, $false, $true, 1, $false, "474×8=", 2)
$d.Content.Find.Execute("971×9=", $true, $false, $false, $false, $false, $true, 1, $false, "139×6=", 2)
$d.Content.Find.Execute("591×7=", $true, $false, $false, $false, $false, $true, 1, $false, "922×8=", 2)
$d.Content.Find.Execute("502×9=", $true, $false, $false, $false, $false, $true, 1, $false, "816×7=", 2)
$d.Content.Find.Execute("780×6=", $true, $false, $false, $false, $false, $true, 1, $false, "139×9=", 2)
$d.Content.Find.Execute("209×8=", $true, $false, $false, $false, $false, $true, 1, $false, "359×4=", 2)
$d.Content.Find.Execute("579×7=", $true, $false, $false, $false, $false, $true, 1, $false, "593×9=", 2)
$d.Content.Find.Execute("175×6=", $true, $false, $false, $false, $false, $true, 1, $false, "572×9=", 2)
$d.Content.Find.Execute("499×7=", $true, $false, $false, $false, $false, $true, 1, $false, "597×9=", 2)
$d.Content.Find.Execute("492×2=", $true, $false, $false, $false, $false, $true, 1, $false, "712×2=", 2)
$d.Content.Find.Execute("800×8=", $true, $false, $false, $false, $false, $true, 1, $false, "128×8=", 2)
